$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 17: add "ALL (458)" label and the new "% AS coverage" result string
$ws.Range("A17").Value = "ALL (458)"
$ws.Range("B17").Value = "5888(with% Ases covered part)"

# Match the formatting already used by neighbouring cells (D16/D17 style)
$ws.Range("D17").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122)

# Extend the chart's plotted series to include the new data point in row 17
$chart = $ws.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "=Sheet1!`$A`$2:`$A`$17"
$series.Values = "=Sheet1!`$B`$2:`$B`$17"
